$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1. Every existing row (1..48) shifts down by one
# (old row 1 -> new row 2, ..., old row 48 -> new row 49).
$ws.Rows.Item(1).Insert()

# The freshly inserted row 1 is blank/unformatted. Copy the header formatting
# (bold / centered / bordered "s=1" style) from row 2 (the old header row,
# now shifted down) onto the new row 1.
$ws.Range("A2:L2").Copy($ws.Range("A1:L1"))

# Populate the new row 1 with the numeric placeholders 0..11.
$ws.Cells.Item(1, 1).Value = 0
$ws.Cells.Item(1, 2).Value = 1
$ws.Cells.Item(1, 3).Value = 2
$ws.Cells.Item(1, 4).Value = 3
$ws.Cells.Item(1, 5).Value = 4
$ws.Cells.Item(1, 6).Value = 5
$ws.Cells.Item(1, 7).Value = 6
$ws.Cells.Item(1, 8).Value = 7
$ws.Cells.Item(1, 9).Value = 8
$ws.Cells.Item(1, 10).Value = 9
$ws.Cells.Item(1, 11).Value = 10
$ws.Cells.Item(1, 12).Value = 11

# Row 2 keeps the old header text ("Lg.", "Threading", ...) but loses the
# header styling (back to the plain/default style used by the other data
# rows), and its K/L cells (thread_size / material_surface) end up blank.
$ws.Range("A2:L2").ClearFormats()
$ws.Cells.Item(2, 11).Value = ""
$ws.Cells.Item(2, 12).Value = ""

Write-Output "edit applied"
